$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New s_vals data (regenerated to filter save games)
$data = @{
    2  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;   E = 0.496779210170732;   G = 9.295990156953671 }
    3  = @{ B = 0.3048080303191223; C = 0.3127903958511391; D = 0.8054896365839992;  E = 8.660232485948974;   G = 10.08332054870323 }
    4  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;   E = 0.496779210170732;   G = 9.295990156953671 }
    5  = @{ B = 0.3048080303191223; C = 25707020678.0705;   D = 3.900430680208489;   E = 645.3272768299601;   G = 25707021327.60302 }
    6  = @{ B = 0.127881588408715;  C = 0.3127903958511391; D = 0.8054896365839992;  E = 0.496779210170732;   G = 1.742940831014585 }
    7  = @{ B = 3.230985683306322;  C = 3099.503889238888;  D = 3.900430680208489;   E = 8.660232485948974;   G = 3115.295538088352 }
    8  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;   E = 0.496779210170732;   G = 9.295990156953671 }
    9  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;   E = 0.496779210170732;   G = 9.295990156953671 }
    10 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615;  E = 0.496779210170732;   G = 5.553084769722144 }
    11 = @{ B = 0.3048080303191223; C = 1.667794583268128;  D = 3.900430680208489;   E = 645.3272768299601;   G = 651.2003101237558 }
    12 = @{ B = 1.459612070389937;  C = 0.04240448674262143;D = 0.1575252929769615;  E = 0.496779210170732;   G = 2.156321060280252 }
    13 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;   E = 8.660232485948974;   G = 17.45944343273191 }
    14 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 26.21740644021617;   E = 8.660232485948974;   G = 38.00504557982321 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
